$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44389
$ws.Range("J2").Value = 81
$ws.Range("K2").Value = 2800
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = 2889
$ws.Range("P2").Value = 963
$ws.Range("D3").Value = 44187
$ws.Range("J3").Value = 65
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 3000
$ws.Range("P3").Value = 1000
$ws.Range("D4").Value = 44669
$ws.Range("J4").Value = 92
$ws.Range("K4").Value = 2500
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 2755
$ws.Range("P4").Value = 918
$ws.Range("D5").Value = 44166
$ws.Range("J5").Value = 45
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = 2500
$ws.Range("P5").Value = 833
$ws.Range("D6").Value = 44557
$ws.Range("J6").Value = 104
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = 2260
$ws.Range("P6").Value = 753
$ws.Range("D7").Value = 44165
$ws.Range("J7").Value = 68
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 3000
$ws.Range("P7").Value = 1000
$ws.Range("D8").Value = 44340
$ws.Range("J8").Value = 54
$ws.Range("K8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 3000
$ws.Range("P8").Value = 1000
$ws.Range("D9").Value = 44224
$ws.Range("J9").Value = 67
$ws.Range("K9").Value = 3000
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 3000
$ws.Range("P9").Value = 1000
$ws.Range("D10").Value = 44242
$ws.Range("J10").Value = 95
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 2737
$ws.Range("P10").Value = 912
$ws.Range("D11").Value = 44223
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 2500
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 2781
$ws.Range("P11").Value = 927
$ws.Range("D12").Value = 44390
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("P12").Value = 1000
$ws.Range("D13").Value = 44536
$ws.Range("J13").Value = 125
$ws.Range("K13").Value = 2200
$ws.Range("L13").Value = 2200
$ws.Range("M13").Value = 2200
$ws.Range("P13").Value = 733
$ws.Range("D14").Value = 44537
$ws.Range("J14").Value = 88
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 2200
$ws.Range("M14").Value = 2091
$ws.Range("P14").Value = 697
$ws.Range("D15").Value = 44292
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("M15").Value = 3000
$ws.Range("P15").Value = 1000
$ws.Range("D16").Value = 44260
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 3500
$ws.Range("L16").Value = 3500
$ws.Range("M16").Value = 3500
$ws.Range("P16").Value = 1167
$ws.Range("D17").Value = 44179
$ws.Range("J17").Value = 78
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = 3000
$ws.Range("P17").Value = 1000
$ws.Range("D18").Value = 44225
$ws.Range("J18").Value = 56
$ws.Range("K18").Value = 3000
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = 3000
$ws.Range("P18").Value = 1000
$ws.Range("D19").Value = 44574
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = 3000
$ws.Range("P19").Value = 1000
$ws.Range("D20").Value = 44559
$ws.Range("J20").Value = 68
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 2000
$ws.Range("P20").Value = 667
$ws.Range("D21").Value = 44627
$ws.Range("J21").Value = 78
$ws.Range("K21").Value = 3500
$ws.Range("L21").Value = 3500
$ws.Range("M21").Value = 3500
$ws.Range("P21").Value = 1167
$ws.Range("D22").Value = 44291
$ws.Range("J22").Value = 45
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = 3000
$ws.Range("P22").Value = 1000
$ws.Range("D23").Value = 44243
$ws.Range("J23").Value = 45
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 3000
$ws.Range("M23").Value = 3000
$ws.Range("P23").Value = 1000
$ws.Range("D24").Value = 44222
$ws.Range("J24").Value = 45
$ws.Range("K24").Value = 3000
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = 3000
$ws.Range("P24").Value = 1000
$ws.Range("D25").Value = 44221
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = 2500
$ws.Range("L25").Value = 2500
$ws.Range("M25").Value = 2500
$ws.Range("P25").Value = 833
$ws.Range("D26").Value = 44193
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = 3000
$ws.Range("L26").Value = 3000
$ws.Range("M26").Value = 3000
$ws.Range("P26").Value = 1000
